$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 573 ("アラビア語で" entry) - remaining rows shift up by one.
$ws.Rows(573).Delete()
